$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append: GSW @ MIN and MIN @ GSW, GM1, 2025-05-06
$row16 = @(14, "GSW", "MIN", "away", "2025-05-06", "240:00", 34, 87, 0.391, 18, 42, 0.429, 13, 15, 0.867, 18, 33, 51, 26, 10, 2, 12, 18, 99, 11, 18, 26, 36, 19, "W")
$row17 = @(15, "MIN", "GSW", "home", "2025-05-06", "240:00", 34, 86, 0.395, 5, 29, 0.172, 15, 17, 0.882, 12, 29, 41, 19, 7, 6, 16, 21, 88, -11, 20, 11, 29, 28, "L")

for ($col = 1; $col -le $row16.Length; $col++) {
    if ($col -eq 5) {
        # The DATE column would otherwise be auto-parsed into a date serial
        # number. Enter it as a formula returning the literal text, then
        # collapse it back down to a plain value so it's stored as text,
        # matching how every other DATE cell in this sheet is stored.
        $ws.Cells.Item(16, $col).Formula = '="' + $row16[$col - 1] + '"'
        $ws.Cells.Item(17, $col).Formula = '="' + $row17[$col - 1] + '"'
    }
    else {
        $ws.Cells.Item(16, $col).Value = $row16[$col - 1]
        $ws.Cells.Item(17, $col).Value = $row17[$col - 1]
    }
}

$ws.Range("E16:E17").Copy()
$ws.Range("E16:E17").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Match the formatting of column A's existing "rank" cells (bold, bordered,
# centered) by copying the format from the row above.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$ws.Range("A1").Select()
